$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text type on the Price/Volume columns so values like "584.09" or
# "0.600" are stored as text (matching original inlineStr cells) instead of
# being auto-converted to numbers by Excel.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '60.219.12'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '2.599.21'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '584.09'
$ws.Range("D6").Value = '142.97'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("E11").Value = '  -1.71%  '
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '3.060.02'
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").Value = '24.49'
$ws.Range("E14").Value = '  +4.10%  '
$ws.Range("D15").Value = '60.224.70'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '2.601.07'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '11.32'
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D20").Value = '345.78'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '0.531'
$ws.Range("E23").Value = '  +2.41%  '
$ws.Range("D24").Value = '63.69'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").Value = '0.160'
$ws.Range("E27").Value = '  +3.19%  '
$ws.Range("E28").Value = '  +4.19%  '
$ws.Range("D29").Value = '0.0₃0797'
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("D30").Value = '6.40'
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").Value = '166.69'
$ws.Range("E32").Value = '  +3.16%  '
$ws.Range("D33").Value = '19.41'
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("E34").Value = '  +9.44%  '
$ws.Range("E35").Value = '  +0.62%  '
$ws.Range("D36").Value = '0.979'
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  +2.70%  '
$ws.Range("D38").Value = '38.18'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").Value = '313.46'
$ws.Range("E39").Value = '  +3.40%  '
$ws.Range("E40").Value = '  +1.60%  '
$ws.Range("E41").Value = '  -1.78%  '
$ws.Range("D42").Value = '135.72'
$ws.Range("E42").Value = '  -3.23%  '
$ws.Range("D43").Value = '0.0994'
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").Value = '19.83'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").Value = '0.0549'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("D49").Value = '0.0241'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '19.86'
$ws.Range("E50").Value = '  +1.85%  '
$ws.Range("E51").Value = '  +0.34%  '

# Restore default (unstyled) appearance now that the text values are set,
# so we do not leave a stray custom number format applied to the cells.
$rng.Style = "Normal"
